# TC01_CDS_Filter_Acesses-Controlled.xlsx
#
# The "startup" sheet has three tabs-worth of rows (CasesTab / SamplesTab /
# FilesTab) whose B (query) and C (StatQuery) columns held large, duplicated
# Cypher query strings. Those queries are removed (cleared out), which also
# drops the now-unused shared strings and lets the row heights -- previously
# stretched tall to show the wrapped query text -- collapse back down to the
# sheet's default height. The previous selection (C4) is moved back to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

# Clear the long Cypher query text out of the query/StatQuery columns for
# the CasesTab, SamplesTab and FilesTab rows (rows 2-4), leaving the
# wrap-text style (s="1") on the cells but with no content.
$ws.Range("B2:C4").ClearContents()

# With the tall wrapped text gone, let the rows shrink back to the sheet's
# normal auto height instead of staying pinned at the old tall heights.
$ws.Range("A1:E6").EntireRow.AutoFit()

# Reset the active selection from C4 back to A2.
$ws.Range("A2").Select()
